$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 454.61
$ws.Range("C2").Value = -0.07
$ws.Range("D2").Value = -0.09
$ws.Range("E2").Value = 9.39
$ws.Range("F2").Value = 1.3
$ws.Range("G2").Value = 7.01
$ws.Range("H2").Value = 20.19
$ws.Range("I2").Value = 16.84
$ws.Range("J2").Value = 79.28
$ws.Range("B3").Value = 14258.49
$ws.Range("C3").Value = -0.16
$ws.Range("D3").Value = -0.05
$ws.Range("E3").Value = 11.49
$ws.Range("F3").Value = 1.59
$ws.Range("G3").Value = 7.69
$ws.Range("H3").Value = 36.23
$ws.Range("I3").Value = 29.81
$ws.Range("J3").Value = 94.51
$ws.Range("B4").Value = 35430.42
$ws.Range("C4").Value = 0.04
$ws.Range("D4").Value = 0.45
$ws.Range("E4").Value = 7.6
$ws.Range("F4").Value = 2.04
$ws.Range("G4").Value = 4.94
$ws.Range("H4").Value = 6.89
$ws.Range("I4").Value = 4.66
$ws.Range("J4").Value = 38.73
$ws.Range("B5").Value = 12.98
$ws.Range("C5").Value = 2.29
$ws.Range("D5").Value = 1.01
$ws.Range("E5").Value = -34.28
$ws.Range("F5").Value = -4.35
$ws.Range("G5").Value = -11.1
$ws.Range("H5").Value = -40.1
$ws.Range("I5").Value = -40.7
$ws.Range("J5").Value = -28.17
$ws.Range("B6").Value = 5.24
$ws.Range("C6").Value = -0.38
$ws.Range("D6").Value = -0.38
$ws.Range("E6").Value = -1.32
$ws.Range("G6").Value = 0.38
$ws.Range("H6").Value = 23.0
$ws.Range("I6").Value = 23.29
$ws.Range("J6").Value = 126.84
$ws.Range("B7").Value = 4.22
$ws.Range("C7").Value = -1.63
$ws.Range("D7").Value = -4.95
$ws.Range("E7").Value = -12.08
$ws.Range("F7").Value = -0.47
$ws.Range("G7").Value = 9.9
$ws.Range("H7").Value = 5.5
$ws.Range("I7").Value = 7.65
$ws.Range("J7").Value = 48.07
$ws.Range("B8").Value = 4.27
$ws.Range("C8").Value = -1.61
$ws.Range("D8").Value = -3.39
$ws.Range("E8").Value = -12.5
$ws.Range("F8").Value = 4.4
$ws.Range("G8").Value = 15.72
$ws.Range("H8").Value = 10.05
$ws.Range("I8").Value = 13.87
$ws.Range("J8").Value = 41.86
$ws.Range("B9").Value = 4.45
$ws.Range("C9").Value = -1.55
$ws.Range("D9").Value = -2.2
$ws.Range("E9").Value = -11.53
$ws.Range("F9").Value = 5.95
$ws.Range("G9").Value = 14.69
$ws.Range("H9").Value = 12.09
$ws.Range("I9").Value = 17.11
$ws.Range("J9").Value = 34.44
$ws.Range("C10").Value = 0.0
$ws.Range("D10").Value = -0.44
$ws.Range("E10").Value = 3.46
$ws.Range("F10").Value = -9.86
$ws.Range("G10").Value = -5.88
$ws.Range("I10").Value = -17.95
$ws.Range("J10").Value = -10.4
$ws.Range("B11").Value = 120.46
$ws.Range("C11").Value = 0.32
$ws.Range("D11").Value = -1.54
$ws.Range("E11").Value = 6.32
$ws.Range("F11").Value = 3.06
$ws.Range("G11").Value = 2.44
$ws.Range("H11").Value = 1.75
$ws.Range("I11").Value = 0.79
$ws.Range("J11").Value = 113.73
$ws.Range("B12").Value = 92.27
$ws.Range("C12").Value = 1.81
$ws.Range("D12").Value = 2.29
$ws.Range("E12").Value = 15.68
$ws.Range("F12").Value = 6.17
$ws.Range("G12").Value = 33.76
$ws.Range("H12").Value = 47.96
$ws.Range("I12").Value = 43.34
$ws.Range("J12").Value = 297.89
$ws.Range("B13").Value = 409.82
$ws.Range("C13").Value = 0.2
$ws.Range("D13").Value = -0.12
$ws.Range("E13").Value = 10.04
$ws.Range("F13").Value = -0.54
$ws.Range("G13").Value = 9.97
$ws.Range("H13").Value = 18.58
$ws.Range("I13").Value = 20.16
$ws.Range("J13").Value = 109.55
$ws.Range("B14").Value = 189.37
$ws.Range("C14").Value = -0.54
$ws.Range("D14").Value = -1.01
$ws.Range("E14").Value = 11.35
$ws.Range("F14").Value = 0.93
$ws.Range("G14").Value = 4.93
$ws.Range("H14").Value = 46.56
$ws.Range("I14").Value = 34.89
$ws.Range("J14").Value = 341.94
$ws.Range("B15").Value = 137.31
$ws.Range("C15").Value = 0.68
$ws.Range("D15").Value = 1.02
$ws.Range("E15").Value = 13.35
$ws.Range("F15").Value = 14.61
$ws.Range("G15").Value = 7.87
$ws.Range("H15").Value = 13.08
$ws.Range("I15").Value = 8.29
$ws.Range("J15").Value = 66.32
$ws.Range("B16").Value = 49.57
$ws.Range("C16").Value = 1.49
$ws.Range("D16").Value = 0.04
$ws.Range("E16").Value = 8.59
$ws.Range("F16").Value = -17.93
$ws.Range("G16").Value = -3.62
$ws.Range("H16").Value = 17.41
$ws.Range("I16").Value = 12.35
$ws.Range("J16").Value = 58.62
$ws.Range("B17").Value = 120.58
$ws.Range("C17").Value = -3.11
$ws.Range("D17").Value = -2.17
$ws.Range("E17").Value = 6.93
$ws.Range("F17").Value = 16.03
$ws.Range("G17").Value = 38.44
$ws.Range("H17").Value = 41.58
$ws.Range("I17").Value = 29.25
$ws.Range("B18").Value = 310.76
$ws.Range("C18").Value = -0.06
$ws.Range("D18").Value = -1.6
$ws.Range("E18").Value = 20.22
$ws.Range("F18").Value = -0.61
$ws.Range("G18").Value = -0.45
$ws.Range("H18").Value = 17.11
$ws.Range("I18").Value = 20.08
$ws.Range("J18").Value = 75.61
$ws.Range("B19").Value = 2047.1
$ws.Range("C19").Value = 0.36
$ws.Range("D19").Value = 2.8
$ws.Range("E19").Value = 2.55
$ws.Range("F19").Value = 5.62
$ws.Range("G19").Value = 4.85
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 17.08
$ws.Range("J19").Value = 67.77
$ws.Range("B20").Value = 25.07
$ws.Range("C20").Value = 0.56
$ws.Range("D20").Value = 5.91
$ws.Range("E20").Value = 7.69
$ws.Range("F20").Value = 2.45
$ws.Range("G20").Value = 6.05
$ws.Range("H20").Value = 5.07
$ws.Range("I20").Value = 18.25
$ws.Range("J20").Value = 77.93
$ws.Range("B21").Value = 936.4
$ws.Range("C21").Value = -0.82
$ws.Range("D21").Value = 1.28
$ws.Range("E21").Value = 0.49
$ws.Range("F21").Value = -3.59
$ws.Range("G21").Value = -7.43
$ws.Range("H21").Value = -12.79
$ws.Range("I21").Value = -8.34
$ws.Range("J21").Value = 17.34
$ws.Range("B22").Value = 1026.6
$ws.Range("C22").Value = -2.91
$ws.Range("D22").Value = -2.62
$ws.Range("E22").Value = -9.5
$ws.Range("F22").Value = -15.0
$ws.Range("G22").Value = -25.97
$ws.Range("H22").Value = -42.63
$ws.Range("I22").Value = -43.48
$ws.Range("J22").Value = -12.15
$ws.Range("B23").Value = 3.79
$ws.Range("C23").Value = -0.26
$ws.Range("D23").Value = 0.53
$ws.Range("E23").Value = 3.84
$ws.Range("F23").Value = 0.53
$ws.Range("G23").Value = 1.88
$ws.Range("H23").Value = -0.52
$ws.Range("I23").Value = 4.41
$ws.Range("J23").Value = 36.33
$ws.Range("B24").Value = 77.86
$ws.Range("C24").Value = 1.9
$ws.Range("D24").Value = 0.99
$ws.Range("E24").Value = -5.41
$ws.Range("F24").Value = -6.9
$ws.Range("G24").Value = 8.53
$ws.Range("H24").Value = -2.99
$ws.Range("I24").Value = -0.43
$ws.Range("J24").Value = 52.88
$ws.Range("B25").Value = 2047.1
$ws.Range("C25").Value = 0.36
$ws.Range("D25").Value = 2.8
$ws.Range("E25").Value = 2.55
$ws.Range("F25").Value = 5.62
$ws.Range("G25").Value = 4.85
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 17.08
$ws.Range("J25").Value = 67.77
$ws.Range("B26").Value = 2.8
$ws.Range("C26").Value = 3.32
$ws.Range("D26").Value = -3.45
$ws.Range("E26").Value = -16.42
$ws.Range("F26").Value = 1.08
$ws.Range("G26").Value = 29.03
$ws.Range("H26").Value = -37.36
$ws.Range("I26").Value = -61.33
$ws.Range("J26").Value = -39.26
$ws.Range("B27").Value = 0.67
$ws.Range("C27").Value = 0.01
$ws.Range("D27").Value = 0.02
$ws.Range("E27").Value = 0.05
$ws.Range("F27").Value = 0.03
$ws.Range("G27").Value = 0.02
$ws.Range("I27").Value = 0.0
$ws.Range("B28").Value = 0.74
$ws.Range("G28").Value = 0.0
$ws.Range("D29").Value = 0.01
$ws.Range("G29").Value = 0.03
$ws.Range("J29").Value = -0.03
$ws.Range("D30").Value = 0.01
$ws.Range("E30").Value = 0.02
$ws.Range("F30").Value = -0.01
$ws.Range("J30").Value = -0.23
$ws.Range("C31").Value = 0.01
$ws.Range("D31").Value = 0.02
$ws.Range("E31").Value = 0.06
$ws.Range("F31").Value = 0.03
$ws.Range("G31").Value = 0.02
$ws.Range("H31").Value = -0.03
$ws.Range("I31").Value = -0.0
$ws.Range("J31").Value = -0.1
$ws.Range("C32").Value = 0.0
$ws.Range("D32").Value = 0.0
$ws.Range("I32").Value = -0.06
$ws.Range("B33").Value = 1.27
$ws.Range("C33").Value = 0.01
$ws.Range("E33").Value = 0.05
$ws.Range("F33").Value = 0.0
$ws.Range("I33").Value = 0.06
$ws.Range("C34").Value = 0.01
$ws.Range("D34").Value = 0.01
$ws.Range("E34").Value = 0.08
$ws.Range("F34").Value = 0.05
$ws.Range("G34").Value = 0.05
$ws.Range("H34").Value = 0.01
$ws.Range("I34").Value = 0.02
$ws.Range("J34").Value = -0.12
$ws.Range("D35").Value = 0.01
$ws.Range("E35").Value = 0.03
$ws.Range("F35").Value = 0.0
$ws.Range("J35").Value = 0.14
